$wb = $excel.ActiveWorkbook

# --- Overview sheet: handoff -> handback status for the second file row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: status + handback datetime + error detail for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-26 10:46:44"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet: status + handback datetime + error detail for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-26 10:46:50"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
